# "Casos de Uso" sheet: fill in the missing value for the 7th use case's
# row so the running totals (shared formulas in AA/AD/.../AZ/BA) recompute.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Casos de Uso")

$ws.Range("Z7").Value = 6

# Leave the view focused where the user left off editing.
$ws.Activate()
$ws.Range("E7").Select()
